# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.528.92"
$ws.Range('E2').Value = '  +2.95%  '

$ws.Range('D3').Value = "'1.854.53"
$ws.Range('E3').Value = '  +1.98%  '

$ws.Range('D4').Value = "'1.005"
$ws.Range('E4').Value = '  +0.48%  '

$ws.Range('D5').Value = "'271.53"
$ws.Range('E5').Value = '  -2.89%  '

$ws.Range('D6').Value = "'1.004"
$ws.Range('E6').Value = '  +0.37%  '

$ws.Range('D7').Value = "'0.5223"
$ws.Range('E7').Value = '  +2.12%  '

$ws.Range('D8').Value = "'0.3375"
$ws.Range('E8').Value = '  -4.81%  '

$ws.Range('D9').Value = "'0.06807"
$ws.Range('E9').Value = '  +2.06%  '

$ws.Range('D10').Value = "'19.74"
$ws.Range('E10').Value = '  -1.82%  '

$ws.Range('D11').Value = "'0.7898"
$ws.Range('E11').Value = '  -4.66%  '

$ws.Range('D12').Value = "'0.07704"
$ws.Range('E12').Value = '  -2.56%  '

$ws.Range('D13').Value = "'1.882.67"
$ws.Range('E13').Value = '  +3.50%  '

$ws.Range('D14').Value = "'89.21"
$ws.Range('E14').Value = '  +1.33%  '

$ws.Range('D15').Value = "'5.110"
$ws.Range('E15').Value = '  +0.44%  '

$ws.Range('D16').Value = "'1.006"
$ws.Range('E16').Value = '  +0.57%  '

$ws.Range('D17').Value = "'14.39"
$ws.Range('E17').Value = '  +2.07%  '

$ws.Range('D18').Value = "'1.003"
$ws.Range('E18').Value = '  +0.32%  '

$ws.Range('D19').Value = "'0.000007943"
$ws.Range('E19').Value = '  -1.22%  '

$ws.Range('D20').Value = "'26.554.22"
$ws.Range('E20').Value = '  +2.89%  '

$ws.Range('D21').Value = "'2.102.39"
$ws.Range('E21').Value = '  +3.16%  '

$ws.Range('D22').Value = "'4.695"
$ws.Range('E22').Value = '  -1.34%  '

$ws.Range('D23').Value = "'9.917"
$ws.Range('E23').Value = '  -0.62%  '

$ws.Range('D24').Value = "'6.105"
$ws.Range('E24').Value = '  -0.37%  '

$ws.Range('D25').Value = "'2.336"
$ws.Range('E25').Value = '  +4.22%  '

$ws.Range('D26').Value = "'145.28"
$ws.Range('E26').Value = '  +2.03%  '

$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').Value = "'1.659"
$ws.Range('E27').Value = '  -0.59%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'17.16"
$ws.Range('E28').Value = '  -0.07%  '

$ws.Range('D29').Value = "'112.21"
$ws.Range('E29').Value = '  +2.57%  '

$ws.Range('D30').Value = "'4.270"
$ws.Range('E30').Value = '  -1.50%  '

$ws.Range('D31').Value = "'4.276"
$ws.Range('E31').Value = '  +0.75%  '

$ws.Range('D32').Value = "'0.08870"
$ws.Range('E32').Value = '  +1.13%  '

$ws.Range('D33').Value = "'0.04895"
$ws.Range('E33').Value = '  -0.37%  '

$ws.Range('D34').Value = "'1.146"
$ws.Range('E34').Value = '  +0.42%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = "'2.888"
$ws.Range('E35').Value = '  +0.54%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = "'0.7188"
$ws.Range('E36').Value = '  -1.92%  '

$ws.Range('D37').Value = "'3.215"
$ws.Range('E37').Value = '  +1.78%  '

$ws.Range('D38').Value = "'0.01835"
$ws.Range('E38').Value = '  -1.04%  '

$ws.Range('D39').Value = "'2.295"
$ws.Range('E39').Value = '  -4.02%  '

$ws.Range('D40').Value = "'0.5056"
$ws.Range('E40').Value = '  -2.17%  '

$ws.Range('D41').Value = "'0.9285"
$ws.Range('E41').Value = '  -3.94%  '

$ws.Range('D42').Value = "'114.96"
$ws.Range('E42').Value = '  +3.32%  '

$ws.Range('D43').Value = "'6.121"
$ws.Range('E43').Value = '  -2.03%  '

$ws.Range('D44').Value = "'7.954"
$ws.Range('E44').Value = '  -1.40%  '

$ws.Range('D45').Value = "'1.004"
$ws.Range('E45').Value = '  +0.38%  '

$ws.Range('D46').Value = "'0.4384"
$ws.Range('E46').Value = '  -4.33%  '

$ws.Range('D47').Value = "'0.1319"
$ws.Range('E47').Value = '  -3.81%  '

$ws.Range('D48').Value = "'9.195"
$ws.Range('E48').Value = '  -0.28%  '

$ws.Range('D49').Value = "'35.95"
$ws.Range('E49').Value = '  -1.95%  '

$ws.Range('D50').Value = "'0.05948"
$ws.Range('E50').Value = '  +2.26%  '

$ws.Range('D51').Value = "'1.465"
$ws.Range('E51').Value = '  -2.65%  '
